$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 16: extra components purchased, with wrapped multi-line description ---
$ws.Range("A16").Value = "N-Mosfet 20V`nP-Mosfet 20V`n68n/50V `n47n/50`nfotorezystor LDR07 50kR`nrezystor pomiarowy 0R1`ndławik drutowy 20u`ntranzystor PNP`n"
$ws.Range("A16").WrapText = $true

$ws.Range("C16").Value = 25

# Enter the date as literal text (avoiding Excel's automatic date recognition)
# via a helper cell computed with TRIM, then paste the resulting text value into
# place; this keeps E16 as a plain shared string without picking up any new
# number-format style.
$ws.Range("Y1").Formula = "=TRIM(""12.10.14 "")"
$ws.Range("Y1").Copy()
$ws.Range("E16").PasteSpecial(-4163)
$ws.Range("Y1").ClearContents()

$ws.Rows.Item(16).RowHeight = 134.25

# Selection moved onto the newly added quantity cell
$ws.Range("C16").Select()
